{"js": "// The paragraph \"<id>p091r_2</id>\" is currently split across four runs:\n//   \"<id>\"    (Courier New, color 7f6000, sz 18)\n//   \"p091r_\"  (color 000000)\n//   \"2\"       (no special formatting)\n//   \"</id>\"   (Courier New, color 7f6000, sz 18)\n// followed by a trailing empty run. Collapse the four runs into a single\n// run that keeps the formatting of the first run (\"<id>\") while holding\n// the combined text \"<id>p091r_2</id>\".\n//\n// Re-inserting the exact same text over the located range (with\n// InsertLocation.replace) merges the underlying runs into one run that\n// inherits the formatting of the range's leading run \u2014 exactly mirroring\n// the OOXML diff.\n\nconst body = context.document.body;\n\n// Exact, literal match on the full tagged id string. It is unique in the\n// document (there is a sibling \"<id>fig_p091v_1</id>\" elsewhere that must\n// stay untouched).\nconst results = body.search(\"<id>p091r_2</id>\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for '<id>p091r_2</id>', found \" + results.items.length);\n}\n\nresults.items[0].insertText(\"<id>p091r_2</id>\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The paragraph \"<id>p091r_2</id>\" is currently split across four runs:\n#   <id>      (Courier New, color 7f6000, sz 18)\n#   p091r_    (color 000000)\n#   2         (no special formatting)\n#   </id>     (Courier New, color 7f6000, sz 18)\n# followed by a trailing empty run. Collapse the last three text runs into\n# the first one, so the paragraph ends up as a single run - keeping the\n# formatting of that first run (\"<id>\") - holding the combined text\n# \"<id>p091r_2</id>\", plus the untouched trailing empty run.\n\n$d = $word.ActiveDocument\n\n$target = \"<id>p091r_2</id>\"\n$prefix = \"<id>\"\n$suffix = $target.Substring($prefix.Length)   # \"p091r_2</id>\"\n\n# Locate the unique occurrence of the full tagged id text. (There is a\n# sibling \"<id>fig_p091v_1</id>\" elsewhere in the document that must stay\n# untouched, so we search for the exact, full string.)\n$whole = $d.Content\n$whole.Find.ClearFormatting()\n$foundWhole = $whole.Find.Execute($target)\nif (-not $foundWhole) {\n    throw \"Could not find target text '$target'\"\n}\n$wholeStart = $whole.Start\n$wholeEnd = $whole.End\n\n# Narrow in on the leading \"<id>\" tag inside that match - its run carries\n# the formatting that must survive on the merged run.\n$openTag = $d.Range($wholeStart, $wholeEnd)\n$openTag.Find.ClearFormatting()\n$foundOpen = $openTag.Find.Execute($prefix)\nif (-not $foundOpen) {\n    throw \"Could not find opening tag '$prefix'\"\n}\n$openTagEnd = $openTag.End\n\n# Delete the remaining (differently formatted) run text that follows\n# \"<id>\" ...\n$rest = $d.Range($openTagEnd, $wholeEnd)\n$rest.Delete()\n\n# ... then append it back onto the surviving \"<id>\" run so the whole\n# thing becomes a single run using that run's formatting.\n$survivingRun = $d.Range($wholeStart, $openTagEnd)\n$survivingRun.InsertAfter($suffix)\n"}
